# edit.ps1
# Applies the batch_test FirstFINALnoExtra_iter2 update described in the commit diff:
#  - updates quantities/rates/descriptions for rows 8-16
#  - inserts a new "Add Tender Premium" line item as row 17
#  - shifts the summary rows (Grand Total / Tender Premium / NET PAYABLE) down by one
#  - refreshes the Grand Total / NET PAYABLE amounts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Medium point -> Long point ---
$ws.Range("C8").Value = 16
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F8").Value = 662
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '10592.00'

# --- Row 9: switch item -> Rewiring of plug point item ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = ''
$ws.Range("C9").Value = 75
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.0'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '0.00'

# --- Row 10: socket item -> wiring run item ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'Mtr.'
$ws.Range("C10").Value = 71
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F10").Value = 81
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '5751.00'

# --- Row 11: fan regulator item -> Plate Earthing item ---
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = 'Set'
$ws.Range("C11").Value = 48
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '13.0'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F11").Value = 5733
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '275184.00'

# --- Row 12: power plug point item -> ceiling fan item ---
$ws.Range("C12").Value = 64
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F12").Value = 1890
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '120960.00'

# --- Row 13: wiring run item -> LED batten item ---
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'Each'
$ws.Range("C13").Value = 87
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F13").Value = 492
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '42804.00'

# --- Row 14: MCB description -> Single pole MCB ---
$ws.Range("C14").Value = 94
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'

# --- Row 15: Single pole MCB -> Double pole MCB ---
$ws.Range("C15").Value = 96
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'

# --- Row 16: 8 Way (8+2) -> 50/63 A rating ---
$ws.Range("C16").Value = 10
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = ' 50/63 A rating'
$ws.Range("F16").Value = 900
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '9000.00'

# --- Insert a new row at 18: pushes the old Grand Total / Tender Premium / NET PAYABLE rows down by one ---
$ws.Rows.Item(18).Insert()

# --- Row 17 (was blank) now holds the new "Add Tender Premium" line item ---
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '%'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 6
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '37'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = 'Add Tender Premium '
$ws.Range("F17").Value = 0
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '0.00'
$ws.Range("H17").Value = 0

# --- Row 18 is the freshly inserted blank row; clear any inherited values beyond column A ---
$ws.Range("B18:I18").ClearContents()

# --- Row 19 (previously row 18): Grand Total Rs., refreshed amount ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '464291.00'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '464291.00'

# --- Row 20 (previously row 19): Tender Premium @ 0% (amounts unchanged) ---

# --- Row 21 (previously row 20): NET PAYABLE AMOUNT Rs., refreshed amount ---
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '464291.00'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '464291.00'

Write-Output "Applied batch_test update: rows 8-16 refreshed, Add Tender Premium row inserted, totals recalculated."
